$wb = $excel.ActiveWorkbook

# --- "cost calculation" sheet ------------------------------------------------
$wsCost = $wb.Worksheets.Item("cost calculation")

# Roster date moved forward one day (warband played another battle).
$wsCost.Range("F2").Value = 43651

# Lost a "forest goblins" henchman model (18 -> 17).
$wsCost.Range("E14").Value = 17

# Opponent / min-foe-warband-rating base value updated.
$wsCost.Range("K13").Value = 192

$wsCost.Range("J9").Select() | Out-Null

# --- "Underdog Bonus" sheet --------------------------------------------------
$wsUnder = $wb.Worksheets.Item("Underdog Bonus")
$wsUnder.Range("D7").Select() | Out-Null

# --- "Characteristic analysis" sheet -----------------------------------------
$wsChar = $wb.Worksheets.Item("Characteristic analysis")
$wsChar.Rows.Item(31).RowHeight = 12.8
$wsChar.Range("M31").Select() | Out-Null

# Leave "cost calculation" as the active sheet/tab, matching the saved file.
$wsCost.Activate() | Out-Null
